# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-10 06:23:04
# Reorders the "G" column (attendance markers) so that a trailing "System" entry
# in the comma-separated list is moved to the front of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    2   = "System, system, backup@backdoor.com"
    3   = "System, dnasr281@gmail.com"
    4   = "System, backup@backdoor.com"
    5   = "System, backup@backdoor.com"
    6   = "System, dnasr281@gmail.com"
    7   = "System, admin@admin.com"
    8   = "System, backup@backdoor.com"
    28  = "System, system, backup@backdoor.com"
    29  = "System, dnasr281@gmail.com"
    30  = "System, backup@backdoor.com"
    31  = "System, backup@backdoor.com"
    32  = "System, dnasr281@gmail.com"
    33  = "System, admin@admin.com"
    34  = "System, backup@backdoor.com"
    54  = "System, system, backup@backdoor.com"
    55  = "System, dnasr281@gmail.com"
    56  = "System, backup@backdoor.com"
    57  = "System, backup@backdoor.com"
    58  = "System, dnasr281@gmail.com"
    59  = "System, admin@admin.com"
    60  = "System, backup@backdoor.com"
    80  = "System, backup@backdoor.com"
    81  = "System, backup@backdoor.com"
    82  = "System, backup@backdoor.com"
    106 = "System, backup@backdoor.com"
    107 = "System, backup@backdoor.com"
    108 = "System, backup@backdoor.com"
    132 = "System, backup@backdoor.com"
    133 = "System, backup@backdoor.com"
    134 = "System, backup@backdoor.com"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
